# Autogenerated on Sun Feb 01 2015 22:24:41 GMT-0500 (Eastern Standard Time)
#
# Renames the "Data" sheet to "Summary" and inserts a new "Source Type"
# header row plus a full OECD source citation block at the bottom of the
# sheet, pushing the existing MSME participation table down by 4 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the sheet ---------------------------------------------------
$ws.Name = "Summary"

# --- Make room for the new "Source Type" row -----------------------------
# Shifts the old rows 5-8 (Micro/SMEs/MSMEs table + source line) down to
# rows 9-12, leaving row 7 free for the new header.
$ws.Range("A5:D8").Insert()

# --- New row 7: Source Type header (bold + underline) --------------------
$ws.Range("A7").Value = "Source Type: Statistical Institution (Most Widely Used)"
$ws.Range("A7").Font.Name = "Calibri"
$ws.Range("A7").Font.Size = 11
$ws.Range("A7").Font.Bold = $true
$ws.Range("A7").Font.Underline = $true

# --- Re-apply formatting lost by the row insert on the rest of the sheet -

# Row 1: Gabon (large "name" style)
$ws.Range("A1").Font.Name = "Calibri"
$ws.Range("A1").Font.Size = 18
$ws.Range("A1").Font.Bold = $false
$ws.Range("A1").Font.Italic = $false
$ws.Range("A1").Font.Underline = $false

# Row 3: MSME Participation on the Economy (bold "title" style)
$ws.Range("A3").Font.Name = "Calibri"
$ws.Range("A3").Font.Size = 11
$ws.Range("A3").Font.Bold = $true
$ws.Range("A3").Font.Italic = $false
$ws.Range("A3").Font.Underline = $false

# Row 9: Micro / SMEs / MSMEs (bold "title" style)
$ws.Range("B9:D9").Font.Name = "Calibri"
$ws.Range("B9:D9").Font.Size = 11
$ws.Range("B9:D9").Font.Bold = $true
$ws.Range("B9:D9").Font.Italic = $false
$ws.Range("B9:D9").Font.Underline = $false

# Row 10: Enterprises (absolute #) row
$ws.Range("A10").Font.Name = "Calibri"
$ws.Range("A10").Font.Size = 11
$ws.Range("A10").Font.Bold = $true
$ws.Range("A10").Font.Italic = $false
$ws.Range("A10").Font.Underline = $false

$ws.Range("C10:D10").Font.Name = "Calibri"
$ws.Range("C10:D10").Font.Size = 11
$ws.Range("C10:D10").Font.Bold = $false
$ws.Range("C10:D10").Font.Italic = $false
$ws.Range("C10:D10").Font.Underline = $false

# Row 11: Enterprises density (per 1000 people) row
$ws.Range("A11").Font.Name = "Calibri"
$ws.Range("A11").Font.Size = 11
$ws.Range("A11").Font.Bold = $true
$ws.Range("A11").Font.Italic = $false
$ws.Range("A11").Font.Underline = $false

$ws.Range("C11:D11").Font.Name = "Calibri"
$ws.Range("C11:D11").Font.Size = 11
$ws.Range("C11:D11").Font.Bold = $false
$ws.Range("C11:D11").Font.Italic = $false
$ws.Range("C11:D11").Font.Underline = $false

# Row 12: Source: OECD, 2005 (italic "source" style)
$ws.Range("A12").Font.Name = "Calibri"
$ws.Range("A12").Font.Size = 11
$ws.Range("A12").Font.Bold = $false
$ws.Range("A12").Font.Italic = $true
$ws.Range("A12").Font.Underline = $false

# --- New rows 20-21: full OECD source citation block ---------------------

# Row 20: OECD (bold "title" style)
$ws.Range("A20").Value = "OECD"
$ws.Range("A20").Font.Name = "Calibri"
$ws.Range("A20").Font.Size = 11
$ws.Range("A20").Font.Bold = $true
$ws.Range("A20").Font.Italic = $false
$ws.Range("A20").Font.Underline = $false

# Row 21: full citation (italic "source" style)
$ws.Range("A21").Value = "Organization for Economic Cooperation and Development (OECD), ""Gabon"", 2005, p. 272. Available at http://www.oecd.org/fr/dev/34883788.pdf"
$ws.Range("A21").Font.Name = "Calibri"
$ws.Range("A21").Font.Size = 11
$ws.Range("A21").Font.Bold = $false
$ws.Range("A21").Font.Italic = $true
$ws.Range("A21").Font.Underline = $false
